$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename clients in the "Nome Cliente" column
$ws.Range("C3").Value = "Newzinho"
$ws.Range("C4").Value = "Matozinhos"

# Move the active selection
$null = $ws.Range("C5").Select()
